$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 326; existing rows 326-376 shift down to 327-377.
$ws.Rows.Item(326).Insert()

# Populate the newly inserted row 326 with the new weekly record.
$ws.Cells.Item(326, 1).Value = 4
$ws.Cells.Item(326, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(326, 3).Value = "Los Lagos"
$ws.Cells.Item(326, 4).Value = 44776
$ws.Cells.Item(326, 5).Value = 10
$ws.Cells.Item(326, 6).Value = 100112008
$ws.Cells.Item(326, 7).Value = "Coliflor"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 100
$ws.Cells.Item(326, 11).Value = 1500
$ws.Cells.Item(326, 12).Value = 1500
$ws.Cells.Item(326, 13).Value = 1500
$ws.Cells.Item(326, 14).Value = "`$/unidad"
$ws.Cells.Item(326, 15).Value = "Región Metropolitana"
$ws.Cells.Item(326, 16).Value = 1500
$ws.Cells.Item(326, 17).Value = 1
$ws.Cells.Item(326, 18).Value = "Hortaliza"
